$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "38.816.79"
$ws.Range("E2").Value = "  +1.72%  "
$ws.Range("D3").Value = "2.095.74"
$ws.Range("E3").Value = "  +0.09%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.07"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.613"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.43"
$ws.Range("E7").Value = "  +1.33%  "
$ws.Range("E9").Value = "  +2.07%  "
$ws.Range("E10").Value = "  -0.08%  "
$ws.Range("E11").Value = "  -0.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.36"
$ws.Range("E12").Value = "  +4.99%  "
$ws.Range("D13").Value = "2.406.36"
$ws.Range("E13").Value = "  +0.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.10"
$ws.Range("E14").Value = "  -0.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.809"
$ws.Range("E15").Value = "  +4.48%  "
$ws.Range("E16").Value = "  -0.29%  "
$ws.Range("D17").Value = "2.095.66"
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("D18").Value = "38.799.66"
$ws.Range("E18").Value = "  +1.77%  "
$ws.Range("E19").Value = "  +2.52%  "
$ws.Range("E20").Value = "  +1.05%  "
$ws.Range("D21").Value = "0.0₃0840"
$ws.Range("E21").Value = "  +0.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.13"
$ws.Range("E22").Value = "  +1.71%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.38"
$ws.Range("E24").Value = "  -0.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.34"
$ws.Range("E25").Value = "  +1.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "171.60"
$ws.Range("E26").Value = "  +0.89%  "
$ws.Range("E27").Value = "  +1.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.138"
$ws.Range("E28").Value = "  +5.60%  "
$ws.Range("E29").Value = "  +3.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.33"
$ws.Range("E30").Value = "  +1.79%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.51"
$ws.Range("E31").Value = "  +4.70%  "
$ws.Range("E32").Value = "  +0.83%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.53"
$ws.Range("E33").Value = "  +2.22%  "
$ws.Range("E34").Value = "  +1.37%  "
$ws.Range("E35").Value = "  +2.43%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.52"
$ws.Range("E36").Value = "  +1.41%  "
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("E38").Value = "  +1.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.13"
$ws.Range("E40").Value = "  +0.67%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0229"
$ws.Range("E41").Value = "  +4.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "101.42"
$ws.Range("E42").Value = "  +1.34%  "
$ws.Range("D43").Value = "1.535.45"
$ws.Range("E43").Value = "  -1.39%  "
$ws.Range("E44").Value = "  -1.34%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.77"
$ws.Range("E45").Value = "  +5.68%  "
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0911"
$ws.Range("E46").Value = "  +0.26%  "
$ws.Range("E47").Value = "  +2.25%  "
$ws.Range("E48").Value = "  -1.31%  "
$ws.Range("E49").Value = "  +1.37%  "
$ws.Range("E50").Value = "  -0.96%  "
$ws.Range("D51").Value = "2.290.93"
$ws.Range("E51").Value = "  -0.03%  "
